# Add "Columbia" targets column to sub_category_scores_v1.xlsx (category_score sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("category_score")

# --- Header cell I1: "Columbia" (copy formatting from the adjacent header H1) ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Columbia"

# --- Numeric target cells in the new column, formatted like the existing blank
#     numeric-style cells in that row (style class "3") ---
$ws.Range("E2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Value = 63

$ws.Range("E3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = 28

$ws.Range("E4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 9

# --- Remaining rows in column I stay blank, but still carry the same blank
#     data-cell style as the rest of the table (copy format only) ---
$ws.Range("B5:B11").Copy()
$ws.Range("I5:I11").PasteSpecial(-4122)

# --- Footer score row: 7.5 like the rest of the row ---
$ws.Range("H12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 7.5

# --- Move the active selection, matching the end-state cursor position ---
$ws.Range("J7").Select()
